# Apply the changes described by the commit:
#  - Move the "metadata" note (old B17 on the data sheet) into its own
#    "metadata" worksheet, stripping the "metadata -> " prefix and
#    prefixing the remaining text with "]" instead.
#  - Fill column A (rows 2-14) of the data sheet with the value 29.
#  - Remove the now-empty/obsolete row 17 from the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Walla Walla - Low Creek")

# 1. Fill A2:A14 with 29.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = 29
}

# 2. Update the note text (was shared by B17): drop the "metadata -> "
#    prefix and use "]" instead.
$oldNote = $ws.Range("B17").Value
$newNote = "]" + "annual redd surveys in Low Creek with experienced and consistent surveyors.  These surveys detected few if any large fluvial fish sized redds, all were resident adults.  "

# 3. Create the new "metadata" worksheet (positioned after the data
#    sheet, matching workbook.xml sheet order) and move the note there.
$metaSheet = $wb.Worksheets.Add($null, $ws)
$metaSheet.Name = "metadata"
$metaSheet.Range("A1").Value = $newNote
$metaSheet.Range("A2").Select() | Out-Null

# 4. Clear out the old note row on the data sheet (row 17), removing it
#    entirely so the sheet's used range shrinks back down to row 14.
$ws.Rows.Item(17).Delete()

# 5. Re-select the data sheet so it stays the active/visible tab, like
#    in the original workbook.
$ws.Select() | Out-Null
